$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows at the top of the data block (before current row 12),
# shifting the existing rows 12-24 down to 14-26.
$ws.Rows("12:13").Insert()

# New row 12: Espárragos, Sin especificar, Primera, $/kilo
$ws.Cells.Item(12, 1).Value = 12
$ws.Cells.Item(12, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(12, 3).Value = "Metropolitana"
$ws.Cells.Item(12, 4).Value = 44460
$ws.Cells.Item(12, 4).Style = $ws.Cells.Item(14, 4).Style
$ws.Cells.Item(12, 4).NumberFormat = $ws.Cells.Item(14, 4).NumberFormat
$ws.Cells.Item(12, 5).Value = 13
$ws.Cells.Item(12, 6).Value = 300000000
$ws.Cells.Item(12, 7).Value = "Espárragos"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 200
$ws.Cells.Item(12, 11).Value = 1800
$ws.Cells.Item(12, 12).Value = 1800
$ws.Cells.Item(12, 13).Value = 1800
$ws.Cells.Item(12, 14).Value = "$/kilo"
$ws.Cells.Item(12, 15).Value = "Provincia de Linares"
$ws.Cells.Item(12, 16).Value = 1800
$ws.Cells.Item(12, 17).Value = 1
$ws.Cells.Item(12, 18).Value = "Hortaliza"

# New row 13: Espárragos, Sin especificar, Segunda, $/kilo
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(13, 3).Value = "Metropolitana"
$ws.Cells.Item(13, 4).Value = 44460
$ws.Cells.Item(13, 4).Style = $ws.Cells.Item(14, 4).Style
$ws.Cells.Item(13, 4).NumberFormat = $ws.Cells.Item(14, 4).NumberFormat
$ws.Cells.Item(13, 5).Value = 13
$ws.Cells.Item(13, 6).Value = 300000000
$ws.Cells.Item(13, 7).Value = "Espárragos"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Segunda"
$ws.Cells.Item(13, 10).Value = 150
$ws.Cells.Item(13, 11).Value = 1500
$ws.Cells.Item(13, 12).Value = 1500
$ws.Cells.Item(13, 13).Value = 1500
$ws.Cells.Item(13, 14).Value = "$/kilo"
$ws.Cells.Item(13, 15).Value = "Provincia de Linares"
$ws.Cells.Item(13, 16).Value = 1500
$ws.Cells.Item(13, 17).Value = 1
$ws.Cells.Item(13, 18).Value = "Hortaliza"
